$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 15:22"

# Update country rows whose ranking/case numbers changed
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 144033
$ws.Range("C8").Value = 309
$ws.Range("D8").Value = 88000
$ws.Range("E8").Value = 51488
$ws.Range("F8").Value = 2922
$ws.Range("G8").Value = 7
$ws.Range("H8").Value = 4545

$ws.Range("A24").Value = "Suecia"
$ws.Range("B24").Value = 14385
$ws.Range("C24").Value = 563
$ws.Range("D24").Value = 550
$ws.Range("E24").Value = 12295
$ws.Range("F24").Value = 450
$ws.Range("G24").Value = 29
$ws.Range("H24").Value = 1540

$ws.Range("A29").Value = "Arabia Saudita"
$ws.Range("B29").Value = 9362
$ws.Range("C29").Value = 1088
$ws.Range("D29").Value = 1398
$ws.Range("E29").Value = 7867
$ws.Range("F29").Value = 78
$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 97

$ws.Range("A30").Value = "Polonia"
$ws.Range("B30").Value = 9082
$ws.Range("C30").Value = 340
$ws.Range("D30").Value = 1040
$ws.Range("E30").Value = 7692
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 350

$ws.Range("A31").Value = "Ecuador"
$ws.Range("B31").Value = 9022
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 1008
$ws.Range("E31").Value = 7558
$ws.Range("F31").Value = 168
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 456

$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 8746
$ws.Range("C32").Value = 328
$ws.Range("D32").Value = 1892
$ws.Range("E32").Value = 6420
$ws.Range("F32").Value = 256
$ws.Range("G32").Value = 13
$ws.Range("H32").Value = 434

$ws.Range("A42").Value = "Serbia"
$ws.Range("B42").Value = 6318
$ws.Range("C42").Value = 324
$ws.Range("D42").Value = 753
$ws.Range("E42").Value = 5443
$ws.Range("F42").Value = 126
$ws.Range("G42").Value = 5
$ws.Range("H42").Value = 122

$ws.Range("A43").Value = "Filipinas"
$ws.Range("B43").Value = 6259
$ws.Range("C43").Value = 172
$ws.Range("D43").Value = 572
$ws.Range("E43").Value = 5278
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 12
$ws.Range("H43").Value = 409

$ws.Range("A64").Value = "Barein"
$ws.Range("B64").Value = 1873
$ws.Range("C64").Value = 100
$ws.Range("D64").Value = 759
$ws.Range("E64").Value = 1107
$ws.Range("F64").Value = 2
$ws.Range("G64").Value = 0
$ws.Range("H64").Value = 7

$ws.Range("A65").Value = "Croacia"
$ws.Range("B65").Value = 1871
$ws.Range("C65").Value = 39
$ws.Range("D65").Value = 709
$ws.Range("E65").Value = 1115
$ws.Range("F65").Value = 23
$ws.Range("G65").Value = 8
$ws.Range("H65").Value = 47

$ws.Range("A68").Value = "Uzbekistan"
$ws.Range("B68").Value = 1543
$ws.Range("C68").Value = 53
$ws.Range("D68").Value = 225
$ws.Range("E68").Value = 1313
$ws.Range("F68").Value = 8
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 5

$ws.Range("A69").Value = "Estonia"
$ws.Range("B69").Value = 1528
$ws.Range("C69").Value = 16
$ws.Range("D69").Value = 164
$ws.Range("E69").Value = 1324
$ws.Range("F69").Value = 10
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 40

$ws.Range("A70").Value = "Irak"
$ws.Range("B70").Value = 1513
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 953
$ws.Range("E70").Value = 478
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 82

$ws.Range("A72").Value = "Azerbaiyan"
$ws.Range("B72").Value = 1398
$ws.Range("C72").Value = 25
$ws.Range("D72").Value = 712
$ws.Range("E72").Value = 667
$ws.Range("F72").Value = 21
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 19

$ws.Range("A115").Value = "Isla de Man"
$ws.Range("B115").Value = 298
$ws.Range("C115").Value = 1
$ws.Range("D115").Value = 192
$ws.Range("E115").Value = 100
$ws.Range("F115").Value = 11
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 6
